# Apply crypto price/volume updates for Sun Feb 19 09:57:14 UTC 2023 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.680.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.698.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3930'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4046'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.18'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08855'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.424'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.109'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001320'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.700.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '99.41'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07043'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.078'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.004'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.73'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.684.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.136'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.60%  '
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.64'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '164.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.802'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +18.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.58'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.144'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09003'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.608'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.069'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02986'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.961'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.03'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.25%  '
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09174'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.463'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7663'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.12%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7173'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.585'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.216'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.003'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.346'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '139.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '90.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.97%  '
